$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "p_pass" column (column Z). This shifts every column to its
# right one place to the left (AA->Z, AB->AA, AC->AB, AD->AC, AE->AD, AF->AE)
# and drops the "p_pass" shared string, so the workbook's dimension goes
# from A1:AF69 to A1:AE69.
$ws.Range("Z1").EntireColumn.Delete()

# The former "Score" column (now AD) held a literal count of TRUE values
# across the five *_pass booleans (p_pass, CFI_pass, GFI_pass, AGFI_pass,
# SRMR_pass). Now that p_pass is gone, recompute it as the count of TRUE
# values across the four remaining *_pass booleans (now in columns Z:AC)
# for every data row.
for ($r = 2; $r -le 69; $r++) {
    $cfiPass   = $ws.Cells.Item($r, 26).Value2   # Z  = CFI_pass
    $gfiPass   = $ws.Cells.Item($r, 27).Value2   # AA = GFI_pass
    $agfiPass  = $ws.Cells.Item($r, 28).Value2   # AB = AGFI_pass
    $srmrPass  = $ws.Cells.Item($r, 29).Value2   # AC = SRMR_pass

    $score = 0
    if ($cfiPass)  { $score = $score + 1 }
    if ($gfiPass)  { $score = $score + 1 }
    if ($agfiPass) { $score = $score + 1 }
    if ($srmrPass) { $score = $score + 1 }

    $ws.Cells.Item($r, 30).Value2 = $score       # AD = Score
}
